# "Generate Report for Handoff"
# Adds two newly-handed-off source files
#   1b336802-f1c0-4209-8bd9-9dbbfa55b339.md
#   d73005fa-e2f0-4523-9198-51bd91932ec4.md
# to the localization-status report: one new row per file on the
# "Overview" sheet (inserted ahead of the existing last row, which is
# pushed down), and one new row per file on each of the "zh-cn" / "de-de"
# detail sheets.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276   # BGR encoding of RGB 6495ED (matches existing "HyperLink" font)
$dateFormat = "yyyy-mm-dd HH:mm:ss"

function Style-Hyperlink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkColor
}

function Style-Date($range) {
    $range.NumberFormat = $dateFormat
}

function Add-Link($ws, $cellRef, $url, $displayText) {
    $range = $ws.Range($cellRef)
    $ws.Hyperlinks.Add($range, $url, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $displayText) | Out-Null
    Style-Hyperlink $range
}

# ---------------------------------------------------------------------
# New source-file rows, keyed the same way for every sheet
# ---------------------------------------------------------------------
$newFiles = @(
    @{
        Guid = "1b336802-f1c0-4209-8bd9-9dbbfa55b339"
        ZhHash = "8c54cafd93eb725680f4b1526cc561de256ca5a7"
        DeHash = "8c54cafd93eb725680f4b1526cc561de256ca5a7"
        MdCommit = "d5f7b6a5e6b0e5c1f3a2b4d6e8f0a1c3b5d7e9f1"
        ZhHandoffCommit = "c1a2b3c4d5e6f7081920a1b2c3d4e5f607182930"
        ZhDatetime = "2016-03-21 10:35:30"
        DeDatetime = "2016-03-21 10:35:36"
        OverviewDatetime = "2016-03-21 10:35:36"
    },
    @{
        Guid = "d73005fa-e2f0-4523-9198-51bd91932ec4"
        ZhHash = "6f3dec9da214b4d05b519db178a92503f4c094c2"
        DeHash = "6f3dec9da214b4d05b519db178a92503f4c094c2"
        MdCommit = "e6f7081920a1b2c3d4e5f607182930c1a2b3c4d5"
        ZhHandoffCommit = "1920a1b2c3d4e5f607182930c1a2b3c4d5e6f708"
        ZhDatetime = "2016-03-21 10:35:30"
        DeDatetime = "2016-03-21 10:35:36"
        OverviewDatetime = "2016-03-21 10:35:36"
    }
)

$statusReady = "Ready for handoff"
$zeroDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Overview sheet: insert the two new rows ahead of the existing last row
# (e011a692-...) which therefore moves from row 5 to row 7.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Rows.Item(5).Insert()
$wsOverview.Rows.Item(5).Insert()

for ($i = 0; $i -lt $newFiles.Count; $i++) {
    $f = $newFiles[$i]
    $row = 5 + $i
    $mdName = "$($f.Guid).md"

    $wsOverview.Cells.Item($row, 1).Value = $mdName
    $wsOverview.Cells.Item($row, 2).Value = $statusReady
    $wsOverview.Cells.Item($row, 3).Value = $statusReady
    $wsOverview.Cells.Item($row, 4).Value = $f.OverviewDatetime
    Style-Date $wsOverview.Range("D$row")

    Add-Link $wsOverview "A$row" "https://github.com/OpenLocalizationTest/oltest/blob/$($f.MdCommit)/e2e/$mdName" $mdName
}

# ---------------------------------------------------------------------
# Detail sheets (zh-cn / de-de)
# ---------------------------------------------------------------------
$langSheets = @(
    @{
        Name = "zh-cn"
        Locale = "zh-cn"
        HashField = "ZhHash"
        DatetimeField = "ZhDatetime"
    },
    @{
        Name = "de-de"
        Locale = "de-de"
        HashField = "DeHash"
        DatetimeField = "DeDatetime"
    }
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)
    $ws.Rows.Item(5).Insert()
    $ws.Rows.Item(5).Insert()

    for ($i = 0; $i -lt $newFiles.Count; $i++) {
        $f = $newFiles[$i]
        $row = 5 + $i
        $mdName = "$($f.Guid).md"
        $hash = $f[$lang.HashField]
        $xlfName = "$($f.Guid).$hash.$($lang.Locale).xlf"
        $datetime = $f[$lang.DatetimeField]

        $ws.Cells.Item($row, 1).Value = $mdName
        $ws.Cells.Item($row, 2).Value = ".md"
        $ws.Cells.Item($row, 3).Value = $statusReady
        $ws.Cells.Item($row, 4).Value = $xlfName
        $ws.Cells.Item($row, 5).Value = $datetime
        Style-Date $ws.Range("E$row")
        $ws.Cells.Item($row, 8).Value = $zeroDate
        Style-Date $ws.Range("H$row")
        $ws.Cells.Item($row, 10).Value = "Include"

        Add-Link $ws "A$row" "https://github.com/OpenLocalizationTest/oltest/blob/$($f.MdCommit)/e2e/$mdName" $mdName

        $handoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$($f.ZhHandoffCommit)/ol-handoff/OpenLocalizationTestOrg/oltest.$($lang.Locale)/ci/ht/$xlfName"
        Add-Link $ws "D$row" $handoffUrl $xlfName
    }
}

Write-Output "done"
